$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.138156
$ws.Range("H2").Value = 0.414468
$ws.Range("I2").Value = 0.0003010053794496939
$ws.Range("J2").Value = 0.0003010053794496939
$ws.Range("M2").Value = 1.01111
$ws.Range("N2").Value = 3.03333
$ws.Range("O2").Value = 0.04063212692754557
$ws.Range("P2").Value = 0.04063212692754556
$ws.Range("Q2").Value = 0.13969091316
$ws.Range("R2").Value = 1.25721821844
$ws.Range("S2").Value = [double]"1.223048878367398E-05"
$ws.Range("T2").Value = [double]"1.223048878367398E-05"

# Row 3
$ws.Range("G3").Value = 0.138156
$ws.Range("H3").Value = 0.414468
$ws.Range("I3").Value = 0.0003010053794496939
$ws.Range("J3").Value = 0.0003010053794496939
$ws.Range("O3").Value = 0.4065982422683317
$ws.Range("P3").Value = 0.4065982422683317
$ws.Range("Q3").Value = 1.397861348804
$ws.Range("R3").Value = 12.580752139236
$ws.Range("S3").Value = 0.0001223882581975578
$ws.Range("T3").Value = 0.0001223882581975578

# Row 4
$ws.Range("G4").Value = 0.138156
$ws.Range("H4").Value = 0.414468
$ws.Range("I4").Value = 0.0003010053794496939
$ws.Range("J4").Value = 0.0003010053794496939
$ws.Range("O4").Value = 0.5527696308041227
$ws.Range("P4").Value = 0.5527696308041226
$ws.Range("Q4").Value = 1.900390167412
$ws.Range("R4").Value = 17.103511506708
$ws.Range("S4").Value = 0.0001663866324684621
$ws.Range("T4").Value = 0.0001663866324684621

# Row 5
$ws.Range("I5").Value = 0.9878623917146768
$ws.Range("J5").Value = 0.9878623917146769
$ws.Range("M5").Value = 1.01111
$ws.Range("N5").Value = 3.03333
$ws.Range("O5").Value = 0.04063212692754557
$ws.Range("P5").Value = 0.04063212692754556
$ws.Range("Q5").Value = 458.4482836397532
$ws.Range("R5").Value = 4126.034552757779
$ws.Range("S5").Value = 0.04013895008709949
$ws.Range("T5").Value = 0.04013895008709948

# Row 6
$ws.Range("I6").Value = 0.9878623917146768
$ws.Range("J6").Value = 0.9878623917146769
$ws.Range("O6").Value = 0.4065982422683317
$ws.Range("P6").Value = 0.4065982422683317
$ws.Range("S6").Value = 0.4016631120741778
$ws.Range("T6").Value = 0.4016631120741778

# Row 7
$ws.Range("I7").Value = 0.9878623917146768
$ws.Range("J7").Value = 0.9878623917146769
$ws.Range("O7").Value = 0.5527696308041227
$ws.Range("P7").Value = 0.5527696308041226
$ws.Range("S7").Value = 0.5460603295533996
$ws.Range("T7").Value = 0.5460603295533994

# Row 8
$ws.Range("G8").Value = 5.432785666666668
$ws.Range("I8").Value = 0.01183660290587349
$ws.Range("J8").Value = 0.01183660290587349
$ws.Range("M8").Value = 1.01111
$ws.Range("N8").Value = 3.03333
$ws.Range("O8").Value = 0.04063212692754557
$ws.Range("P8").Value = 0.04063212692754556
$ws.Range("Q8").Value = 5.493143915423334
$ws.Range("R8").Value = 49.43829523881001
$ws.Range("S8").Value = 0.0004809463516624065
$ws.Range("T8").Value = 0.0004809463516624064

# Row 9
$ws.Range("G9").Value = 5.432785666666668
$ws.Range("I9").Value = 0.01183660290587349
$ws.Range("J9").Value = 0.01183660290587349
$ws.Range("O9").Value = 0.4065982422683317
$ws.Range("P9").Value = 0.4065982422683317
$ws.Range("Q9").Value = 54.96888372397657
$ws.Range("R9").Value = 494.7199535157891
$ws.Range("S9").Value = 0.00481274193595639
$ws.Range("T9").Value = 0.00481274193595639

# Row 10
$ws.Range("G10").Value = 5.432785666666668
$ws.Range("I10").Value = 0.01183660290587349
$ws.Range("J10").Value = 0.01183660290587349
$ws.Range("O10").Value = 0.5527696308041227
$ws.Range("P10").Value = 0.5527696308041226
$ws.Range("Q10").Value = 74.73010555162413
$ws.Range("R10").Value = 672.5709499646172
$ws.Range("S10").Value = 0.006542914618254697
$ws.Range("T10").Value = 0.006542914618254695
